{"js": "// Replace the 25 division problems in the worksheet table with new values,\n// matching the exact text runs in document order (5 content rows x 5 columns,\n// with blank spacer rows in between). Only the <w:t> text content is changed;\n// run/paragraph formatting (fonts, size, alignment) is left untouched.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Content lives in rows 0, 4, 8, 12, 16 (rows 1-3, 5-7, ... are blank spacers).\nconst contentRows = [0, 4, 8, 12, 16];\n\n// Replacement values, in left-to-right / top-to-bottom document order,\n// 5 per content row.\nconst replacements = [\n  [\"12\u00f72=\", \"35\u00f75=\"],\n  [\"46\u00f75=\", \"11\u00f75=\"],\n  [\"58\u00f72=\", \"22\u00f76=\"],\n  [\"69\u00f79=\", \"53\u00f76=\"],\n  [\"71\u00f72=\", \"69\u00f78=\"],\n  [\"64\u00f74=\", \"87\u00f72=\"],\n  [\"40\u00f75=\", \"43\u00f74=\"],\n  [\"57\u00f77=\", \"89\u00f76=\"],\n  [\"11\u00f75=\", \"14\u00f76=\"],\n  [\"87\u00f74=\", \"67\u00f72=\"],\n  [\"51\u00f75=\", \"23\u00f76=\"],\n  [\"77\u00f75=\", \"79\u00f77=\"],\n  [\"48\u00f72=\", \"44\u00f78=\"],\n  [\"99\u00f77=\", \"98\u00f72=\"],\n  [\"41\u00f78=\", \"74\u00f78=\"],\n  [\"36\u00f74=\", \"45\u00f72=\"],\n  [\"84\u00f79=\", \"58\u00f76=\"],\n  [\"82\u00f75=\", \"43\u00f72=\"],\n  [\"12\u00f77=\", \"75\u00f79=\"],\n  [\"53\u00f77=\", \"68\u00f75=\"],\n  [\"23\u00f78=\", \"10\u00f72=\"],\n  [\"22\u00f72=\", \"82\u00f75=\"],\n  [\"43\u00f73=\", \"40\u00f77=\"],\n  [\"58\u00f77=\", \"97\u00f77=\"],\n  [\"32\u00f79=\", \"51\u00f75=\"],\n];\n\n// Gather the target cells first and verify their current text matches what\n// the diff expects before touching anything, so we never clobber the wrong\n// cell if the table shape is not what we assume. (cell.body.text is the\n// clean paragraph text; a range over the cell body additionally reports a\n// trailing tab/cell-mark, so we use body.text purely for the sanity check.)\nconst cells = [];\nlet idx = 0;\nfor (const rowIndex of contentRows) {\n  for (let col = 0; col < 5; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.body.load(\"text\");\n    cells.push(cell);\n    idx++;\n  }\n}\nawait context.sync();\n\nidx = 0;\nfor (const rowIndex of contentRows) {\n  for (let col = 0; col < 5; col++) {\n    const [oldText, newText] = replacements[idx];\n    const cell = cells[idx];\n    idx++;\n\n    if (cell.body.text !== oldText) {\n      throw new Error(\n        \"Unexpected cell text at row \" + rowIndex + \", col \" + col +\n        \": expected '\" + oldText + \"' but found '\" + cell.body.text + \"'\"\n      );\n    }\n\n    // Overwrite just the text; formatting of the run/paragraph is preserved\n    // because we replace the range's text in place rather than clearing it.\n    const cellRange = cell.body.getRange();\n    cellRange.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division problems in the worksheet table with new values,\n# matching the exact text runs in document order (5 content rows x 5 columns,\n# with blank spacer rows in between). Setting Cell.Range.Text preserves the\n# existing run/paragraph formatting (fonts, size, alignment) since Word keeps\n# the first run's rPr when the cell range text is replaced in place.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Content lives in (1-indexed) rows 1, 5, 9, 13, 17; the rows in between are\n# blank spacer rows.\n$contentRows = @(1, 5, 9, 13, 17)\n\n# Replacement values, in left-to-right / top-to-bottom document order,\n# 5 per content row: old text (sanity reference) -> new text.\n$replacements = @(\n  @(\"12\u00f72=\", \"35\u00f75=\"),\n  @(\"46\u00f75=\", \"11\u00f75=\"),\n  @(\"58\u00f72=\", \"22\u00f76=\"),\n  @(\"69\u00f79=\", \"53\u00f76=\"),\n  @(\"71\u00f72=\", \"69\u00f78=\"),\n  @(\"64\u00f74=\", \"87\u00f72=\"),\n  @(\"40\u00f75=\", \"43\u00f74=\"),\n  @(\"57\u00f77=\", \"89\u00f76=\"),\n  @(\"11\u00f75=\", \"14\u00f76=\"),\n  @(\"87\u00f74=\", \"67\u00f72=\"),\n  @(\"51\u00f75=\", \"23\u00f76=\"),\n  @(\"77\u00f75=\", \"79\u00f77=\"),\n  @(\"48\u00f72=\", \"44\u00f78=\"),\n  @(\"99\u00f77=\", \"98\u00f72=\"),\n  @(\"41\u00f78=\", \"74\u00f78=\"),\n  @(\"36\u00f74=\", \"45\u00f72=\"),\n  @(\"84\u00f79=\", \"58\u00f76=\"),\n  @(\"82\u00f75=\", \"43\u00f72=\"),\n  @(\"12\u00f77=\", \"75\u00f79=\"),\n  @(\"53\u00f77=\", \"68\u00f75=\"),\n  @(\"23\u00f78=\", \"10\u00f72=\"),\n  @(\"22\u00f72=\", \"82\u00f75=\"),\n  @(\"43\u00f73=\", \"40\u00f77=\"),\n  @(\"58\u00f77=\", \"97\u00f77=\"),\n  @(\"32\u00f79=\", \"51\u00f75=\")\n)\n\n$idx = 0\nforeach ($row in $contentRows) {\n  for ($col = 1; $col -le 5; $col++) {\n    $pair = $replacements[$idx]\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $idx++\n\n    $cell = $t.Cell($row, $col)\n\n    # Cell.Range.Text includes the trailing cell-end mark (CR + BEL); strip\n    # it before comparing so we can confirm we are editing the expected\n    # cell before touching anything.\n    $actual = $cell.Range.Text.TrimEnd([char]7).TrimEnd([char]13)\n    if ($actual -ne $oldText) {\n      throw \"Unexpected cell text at row $row, col $col`: expected '$oldText' but found '$actual'\"\n    }\n\n    # Overwrite just the text; Word preserves the existing run formatting\n    # (rFonts/sz) and paragraph formatting (jc) of the cell when the range\n    # text is replaced in place like this.\n    $cell.Range.Text = $newText\n  }\n}\n"}
